# "Generate Report for Handoff"
# Updates the localization-status report: the zh-cn / de-de items have moved
# from "In Translation" to "Ready for handoff", so the Status columns and the
# relevant Latest Handoff/Xliff-generation timestamps are refreshed, and the
# (now wider) Status columns are resized to fit the new text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Status: "In Translation" -> "Ready for handoff" -----------------------
$overview.Range("E2").Value = "Ready for handoff"   # zh-cn status
$overview.Range("F2").Value = "Ready for handoff"   # de-de status
$zhcn.Range("C2").Value     = "Ready for handoff"   # Status
$dede.Range("C2").Value     = "Ready for handoff"   # Status

# --- Refreshed handoff / xliff-generation timestamps ------------------------
$zhcn.Range("H2").Value     = "2016-08-28 20:39:14" # Latest Handoff Datetime (zh-cn)
$dede.Range("H2").Value     = "2016-08-28 20:39:19" # Latest Handoff Datetime (de-de)
$overview.Range("G2").Value = "2016-08-28 20:39:19" # Latest HO Xliff Generate Date

# --- Resize the Status columns to fit the new, longer text ------------------
$overview.Columns.Item(5).ColumnWidth = 16.3   # column E (zh-cn status)
$overview.Columns.Item(6).ColumnWidth = 16.3   # column F (de-de status)
$zhcn.Columns.Item(3).ColumnWidth     = 16.3   # column C (Status)
$dede.Columns.Item(3).ColumnWidth     = 16.3   # column C (Status)
